$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume/Coin/Link columns keep their numeric-looking values
# as literal text (matching the original inlineStr cells) instead of being
# auto-converted to numbers by Excel.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '69.115.68'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '3.758.79'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '602.63'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Value = '166.65'
$ws.Range('E6').Value = '  -1.20%  '
$ws.Range('D7').Value = '3.756.54'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.538'
$ws.Range('E9').Value = '  +0.51%  '
$ws.Range('D10').Value = '0.171'
$ws.Range('E10').Value = '  +4.36%  '
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('D12').Value = '0.459'
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('D13').Value = '37.72'
$ws.Range('E13').Value = '  -1.28%  '
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('D15').Value = '4.386.43'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '3.762.92'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').Value = '69.114.62'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').Value = '7.39'
$ws.Range('E18').Value = '  +1.65%  '
$ws.Range('D19').Value = '17.71'
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('D21').Value = '11.33'
$ws.Range('E21').Value = '  +4.12%  '
$ws.Range('D22').Value = '491.29'
$ws.Range('E22').Value = '  -0.83%  '
$ws.Range('E23').Value = '  -0.25%  '
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').Value = '0.0000149'
$ws.Range('E24').Value = '  -0.94%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '84.70'
$ws.Range('E25').Value = '  -0.97%  '
$ws.Range('E26').Value = '  -2.25%  '
$ws.Range('E27').Value = '  -0.42%  '
$ws.Range('D28').Value = '10.07'
$ws.Range('E28').Value = '  -1.26%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  -0.48%  '
$ws.Range('E31').Value = '  +2.27%  '
$ws.Range('D32').Value = '2.42'
$ws.Range('E32').Value = '  -3.96%  '
$ws.Range('D33').Value = '31.74'
$ws.Range('E33').Value = '  -0.77%  '
$ws.Range('D34').Value = '3.903.78'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').Value = '3.697.82'
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('E36').Value = '  -0.30%  '
$ws.Range('D37').Value = '0.139'
$ws.Range('E37').Value = '  +5.98%  '
$ws.Range('D38').Value = '5.94'
$ws.Range('E38').Value = '  +1.55%  '
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('E41').Value = '  +8.35%  '
$ws.Range('D42').Value = '0.325'
$ws.Range('E42').Value = '  +0.65%  '
$ws.Range('D43').Value = '428.39'
$ws.Range('E43').Value = '  -3.47%  '
$ws.Range('E44').Value = '  +1.13%  '
$ws.Range('D45').Value = '48.59'
$ws.Range('E45').Value = '  -0.81%  '
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').Value = '40.39'
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('D49').Value = '142.69'
$ws.Range('E49').Value = '  +0.70%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.810.07'
$ws.Range('E50').Value = '  -0.74%  '
$ws.Range('B51').Value = 'ONDO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D51').Value = '1.30'
$ws.Range('E51').Value = '  +8.69%  '
